$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 with revised C (score) and D (importances) values
$ws.Range("C2").Value = 0.6389872355837201
$ws.Range("D2").Value = 0

$ws.Range("C3").Value = -0.7868771669056187
$ws.Range("D3").Value = 61.35875894829685

$ws.Range("C4").Value = 0.4765274912498819
$ws.Range("D4").Value = 3.539790619544031

$ws.Range("C5").Value = 0.5981626139410883
$ws.Range("D5").Value = 2.551261415002377

$ws.Range("C6").Value = 0.5870731680754302
$ws.Range("D6").Value = 4.19520554511559

$ws.Range("C7").Value = 0.5620547403671418
$ws.Range("D7").Value = 2.647492674317471

$ws.Range("C8").Value = 0.5995260410612524
$ws.Range("D8").Value = 2.08734404593266

$ws.Range("C9").Value = 0.4937239777967523
$ws.Range("D9").Value = 4.06942850673332

$ws.Range("C10").Value = 0.4462461032289393
$ws.Range("D10").Value = 7.880109977251994

# Add new row 11
$ws.Range("A11").Value = 9
$ws.Range("A11").Font.Bold = $true
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("A11").Borders.LineStyle = 1
$ws.Range("A11").Borders.Weight = 2

$ws.Range("B11").Value = "NP_incubation Concentration (mg/mL)"

$ws.Range("C11").Value = 0.5851445753972768
$ws.Range("D11").Value = 11.67060826780572
